# Generate Report for Handoff
# Updates the localization-status workbook to reflect a fresh handoff report:
#   - Status cells move from "Handed back: in sync with en-US" to "Ready for handoff"
#   - The associated generation/handoff timestamps are refreshed
#   - The now-shorter "Status" column narrows (was sized for the long "Handed back..."
#     text, now sized for "Ready for handoff")

$wb = $excel.ActiveWorkbook

$wsOverview = $wb.Worksheets.Item("Overview")
$wsZhCn     = $wb.Worksheets.Item("zh-cn")
$wsDeDe     = $wb.Worksheets.Item("de-de")

# --- Overview sheet -------------------------------------------------------
# E2 (zh-cn status) and F2 (de-de status): "Handed back..." -> "Ready for handoff"
$wsOverview.Range("E2").Value = "Ready for handoff"
$wsOverview.Range("F2").Value = "Ready for handoff"
# G2: Latest HO Xliff Generate Date refreshed
$wsOverview.Range("G2").Value = "2016-08-31 08:25:43"

# Columns E and F were sized for the long status text; narrow them to fit
# the shorter "Ready for handoff" label.
$wsOverview.Columns.Item(5).ColumnWidth = 16.333333333333332
$wsOverview.Columns.Item(6).ColumnWidth = 16.333333333333332

# --- zh-cn sheet ------------------------------------------------------------
$wsZhCn.Range("C2").Value = "Ready for handoff"
$wsZhCn.Range("H2").Value = "2016-08-31 08:25:31"
$wsZhCn.Columns.Item(3).ColumnWidth = 16.333333333333332

# --- de-de sheet ------------------------------------------------------------
$wsDeDe.Range("C2").Value = "Ready for handoff"
$wsDeDe.Columns.Item(3).ColumnWidth = 16.333333333333332
